$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-03 Wednesday", "2024-04-04 Thursday"),
    @("624×8=", "501×3="),
    @("161×9=", "941×3="),
    @("628×7=", "811×9="),
    @("341×8=", "657×3="),
    @("980×9=", "551×2="),
    @("945×5=", "702×5="),
    @("275×8=", "204×8="),
    @("881×5=", "300×6="),
    @("612×4=", "575×6="),
    @("403×9=", "862×5="),
    @("914×9=", "888×8="),
    @("430×7=", "228×8="),
    @("293×4=", "195×6="),
    @("743×5=", "953×4="),
    @("431×6=", "869×8="),
    @("333×9=", "787×6="),
    @("918×9=", "203×3="),
    @("925×8=", "264×6="),
    @("690×2=", "279×7="),
    @("791×8=", "752×8="),
    @("784×3=", "898×8="),
    @("548×5=", "818×5="),
    @("699×4=", "348×9="),
    @("193×3=", "118×2="),
    @("630×2=", "846×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
